# Added pdf output functionality
# The Output column (column B) on Sheet1 is updated to use the new
# A/B/C/D/E output-channel labels instead of the old D/f/I/L/C labels.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$ws1.Range("B2").Value = "A"
$ws1.Range("B3").Value = "B"
$ws1.Range("B4").Value = "C"
$ws1.Range("B5").Value = "D"
$ws1.Range("B6").Value = "E"

# Move the selection on Sheet1 to B6 (matches the saved selection state),
# then restore Sheet2 as the active sheet/tab, since Sheet2 was the tab
# shown when the workbook was originally saved.
[void]$ws1.Activate()
[void]$ws1.Range("B6").Select()
[void]$ws2.Activate()
